$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.505614041169197
$ws.Cells.Item(2, 3).Value = 1.65323645889881
$ws.Cells.Item(2, 4).Value = 3.082599426703578
$ws.Cells.Item(2, 5).Value = 6.48142807727062
$ws.Cells.Item(2, 7).Value = 12.7228780040422
$ws.Cells.Item(3, 2).Value = 3.182878228561681
$ws.Cells.Item(3, 3).Value = 1.65323645889881
$ws.Cells.Item(3, 4).Value = 0.7127328510149897
$ws.Cells.Item(3, 5).Value = 0.4998867070740569
$ws.Cells.Item(3, 7).Value = 6.048734245549538
$ws.Cells.Item(4, 2).Value = 1.505614041169197
$ws.Cells.Item(4, 3).Value = 1.65323645889881
$ws.Cells.Item(4, 4).Value = 0.7127328510149897
$ws.Cells.Item(4, 5).Value = 0.4998867070740569
$ws.Cells.Item(4, 7).Value = 4.371470058157054
$ws.Cells.Item(5, 2).Value = 3.182878228561681
$ws.Cells.Item(5, 3).Value = 1.65323645889881
$ws.Cells.Item(5, 4).Value = 3.082599426703578
$ws.Cells.Item(5, 5).Value = 0.4998867070740569
$ws.Cells.Item(5, 7).Value = 8.418600821238126
$ws.Cells.Item(6, 2).Value = 0.7287194209349384
$ws.Cells.Item(6, 3).Value = 1.65323645889881
$ws.Cells.Item(6, 4).Value = 0.7127328510149897
$ws.Cells.Item(6, 5).Value = 0.4998867070740569
$ws.Cells.Item(6, 7).Value = 3.594575437922795
$ws.Cells.Item(7, 2).Value = 0.001754667048134761
$ws.Cells.Item(7, 3).Value = 0.0001537489499301437
$ws.Cells.Item(7, 4).Value = 0.7127328510149897
$ws.Cells.Item(7, 5).Value = 0.4998867070740569
$ws.Cells.Item(7, 7).Value = 1.214527974087112
$ws.Cells.Item(8, 2).Value = 3.182878228561681
$ws.Cells.Item(8, 3).Value = 1.65323645889881
$ws.Cells.Item(8, 4).Value = 3.082599426703578
$ws.Cells.Item(8, 5).Value = 0.4998867070740569
$ws.Cells.Item(8, 7).Value = 8.418600821238126
$ws.Cells.Item(9, 2).Value = 3.182878228561681
$ws.Cells.Item(9, 3).Value = 1.65323645889881
$ws.Cells.Item(9, 4).Value = 0.7127328510149897
$ws.Cells.Item(9, 5).Value = 0.4998867070740569
$ws.Cells.Item(9, 7).Value = 6.048734245549538
$ws.Cells.Item(10, 2).Value = 3.182878228561681
$ws.Cells.Item(10, 3).Value = 1.65323645889881
$ws.Cells.Item(10, 4).Value = 0.7127328510149897
$ws.Cells.Item(10, 5).Value = 6.48142807727062
$ws.Cells.Item(10, 7).Value = 12.0302756157461
$ws.Cells.Item(11, 2).Value = 1.505614041169197
$ws.Cells.Item(11, 3).Value = 1.65323645889881
$ws.Cells.Item(11, 4).Value = 0.7127328510149897
$ws.Cells.Item(11, 5).Value = 0.4998867070740569
$ws.Cells.Item(11, 7).Value = 4.371470058157054
$ws.Cells.Item(12, 2).Value = 1.505614041169197
$ws.Cells.Item(12, 3).Value = 1.65323645889881
$ws.Cells.Item(12, 4).Value = 3.082599426703578
$ws.Cells.Item(12, 5).Value = 0.4998867070740569
$ws.Cells.Item(12, 7).Value = 6.741336633845642
$ws.Cells.Item(13, 2).Value = 3.182878228561681
$ws.Cells.Item(13, 3).Value = 1.65323645889881
$ws.Cells.Item(13, 4).Value = 0.1529057820181812
$ws.Cells.Item(13, 5).Value = 0.4998867070740569
$ws.Cells.Item(13, 7).Value = 5.488907176552729
$ws.Cells.Item(14, 2).Value = 3.182878228561681
$ws.Cells.Item(14, 3).Value = 1.65323645889881
$ws.Cells.Item(14, 4).Value = 0.1529057820181812
$ws.Cells.Item(14, 5).Value = 0.4998867070740569
$ws.Cells.Item(14, 7).Value = 5.488907176552729
$ws.Cells.Item(15, 2).Value = 3.182878228561681
$ws.Cells.Item(15, 3).Value = 1.65323645889881
$ws.Cells.Item(15, 4).Value = 3.082599426703578
$ws.Cells.Item(15, 5).Value = 0.4998867070740569
$ws.Cells.Item(15, 7).Value = 8.418600821238126
$ws.Cells.Item(16, 2).Value = 3.182878228561681
$ws.Cells.Item(16, 3).Value = 1.65323645889881
$ws.Cells.Item(16, 4).Value = 16.98373111632243
$ws.Cells.Item(16, 5).Value = 0.4998867070740569
$ws.Cells.Item(16, 7).Value = 22.31973251085698
$ws.Cells.Item(17, 2).Value = 3.182878228561681
$ws.Cells.Item(17, 3).Value = 1.65323645889881
$ws.Cells.Item(17, 4).Value = 0.1529057820181812
$ws.Cells.Item(17, 5).Value = 0.4998867070740569
$ws.Cells.Item(17, 7).Value = 5.488907176552729
$ws.Cells.Item(18, 2).Value = 0.1554434735375247
$ws.Cells.Item(18, 3).Value = 0.05231270169004087
$ws.Cells.Item(18, 4).Value = 0.1529057820181812
$ws.Cells.Item(18, 5).Value = 0.4998867070740569
$ws.Cells.Item(18, 7).Value = 0.8605486643198037
$ws.Cells.Item(19, 2).Value = 0.1554434735375247
$ws.Cells.Item(19, 3).Value = 86.29678392075563
$ws.Cells.Item(19, 4).Value = 0.1529057820181812
$ws.Cells.Item(19, 5).Value = 6.48142807727062
$ws.Cells.Item(19, 7).Value = 93.08656125358196
$ws.Cells.Item(20, 2).Value = 1.505614041169197
$ws.Cells.Item(20, 3).Value = 1.65323645889881
$ws.Cells.Item(20, 4).Value = 0.7127328510149897
$ws.Cells.Item(20, 5).Value = 6.48142807727062
$ws.Cells.Item(20, 7).Value = 10.35301142835362
$ws.Cells.Item(21, 2).Value = 1.505614041169197
$ws.Cells.Item(21, 3).Value = 1.65323645889881
$ws.Cells.Item(21, 4).Value = 0.1529057820181812
$ws.Cells.Item(21, 5).Value = 0.4998867070740569
$ws.Cells.Item(21, 7).Value = 3.811642989160245
$ws.Cells.Item(22, 2).Value = 1.505614041169197
$ws.Cells.Item(22, 3).Value = 1.65323645889881
$ws.Cells.Item(22, 4).Value = 0.7127328510149897
$ws.Cells.Item(22, 5).Value = 0.4998867070740569
$ws.Cells.Item(22, 7).Value = 4.371470058157054
$ws.Cells.Item(23, 2).Value = 3.182878228561681
$ws.Cells.Item(23, 3).Value = 1.65323645889881
$ws.Cells.Item(23, 4).Value = 0.7127328510149897
$ws.Cells.Item(23, 5).Value = 0.4998867070740569
$ws.Cells.Item(23, 7).Value = 6.048734245549538
$ws.Cells.Item(24, 2).Value = 3.182878228561681
$ws.Cells.Item(24, 3).Value = 1.65323645889881
$ws.Cells.Item(24, 4).Value = 0.7127328510149897
$ws.Cells.Item(24, 5).Value = 0.4998867070740569
$ws.Cells.Item(24, 7).Value = 6.048734245549538
$ws.Cells.Item(25, 2).Value = 0.3464964993005633
$ws.Cells.Item(25, 3).Value = 0.004309184025731883
$ws.Cells.Item(25, 4).Value = 3.082599426703578
$ws.Cells.Item(25, 5).Value = 0.4998867070740569
$ws.Cells.Item(25, 7).Value = 3.933291817103931
$ws.Cells.Item(26, 2).Value = 3.182878228561681
$ws.Cells.Item(26, 3).Value = 1.65323645889881
$ws.Cells.Item(26, 4).Value = 0.1529057820181812
$ws.Cells.Item(26, 5).Value = 0.4998867070740569
$ws.Cells.Item(26, 7).Value = 5.488907176552729
$ws.Cells.Item(27, 2).Value = 3.182878228561681
$ws.Cells.Item(27, 3).Value = 1.65323645889881
$ws.Cells.Item(27, 4).Value = 0.1529057820181812
$ws.Cells.Item(27, 5).Value = 0.4998867070740569
$ws.Cells.Item(27, 7).Value = 5.488907176552729
$ws.Cells.Item(28, 2).Value = 1.505614041169197
$ws.Cells.Item(28, 3).Value = 1.65323645889881
$ws.Cells.Item(28, 4).Value = 0.7127328510149897
$ws.Cells.Item(28, 5).Value = 0.4998867070740569
$ws.Cells.Item(28, 7).Value = 4.371470058157054
$ws.Cells.Item(29, 2).Value = 1.505614041169197
$ws.Cells.Item(29, 3).Value = 1.65323645889881
$ws.Cells.Item(29, 4).Value = 0.1529057820181812
$ws.Cells.Item(29, 5).Value = 0.4998867070740569
$ws.Cells.Item(29, 7).Value = 3.811642989160245
$ws.Cells.Item(30, 2).Value = 3.182878228561681
$ws.Cells.Item(30, 3).Value = 1.65323645889881
$ws.Cells.Item(30, 4).Value = 0.7127328510149897
$ws.Cells.Item(30, 5).Value = 0.4998867070740569
$ws.Cells.Item(30, 7).Value = 6.048734245549538
$ws.Cells.Item(31, 2).Value = 3.182878228561681
$ws.Cells.Item(31, 3).Value = 1.65323645889881
$ws.Cells.Item(31, 4).Value = 3.082599426703578
$ws.Cells.Item(31, 5).Value = 0.4998867070740569
$ws.Cells.Item(31, 7).Value = 8.418600821238126
$ws.Cells.Item(32, 2).Value = 3.182878228561681
$ws.Cells.Item(32, 3).Value = 1.65323645889881
$ws.Cells.Item(32, 4).Value = 0.7127328510149897
$ws.Cells.Item(32, 5).Value = 0.4998867070740569
$ws.Cells.Item(32, 7).Value = 6.048734245549538
$ws.Cells.Item(33, 2).Value = 0.06328177979961902
$ws.Cells.Item(33, 3).Value = 1.65323645889881
$ws.Cells.Item(33, 4).Value = 0.1529057820181812
$ws.Cells.Item(33, 5).Value = 0.4998867070740569
$ws.Cells.Item(33, 7).Value = 2.369310727790667
$ws.Cells.Item(34, 2).Value = 1.505614041169197
$ws.Cells.Item(34, 3).Value = 1.65323645889881
$ws.Cells.Item(34, 4).Value = 0.7127328510149897
$ws.Cells.Item(34, 5).Value = 6.48142807727062
$ws.Cells.Item(34, 7).Value = 10.35301142835362
$ws.Cells.Item(35, 2).Value = 3.182878228561681
$ws.Cells.Item(35, 3).Value = 1.65323645889881
$ws.Cells.Item(35, 4).Value = 3.082599426703578
$ws.Cells.Item(35, 5).Value = 0.4998867070740569
$ws.Cells.Item(35, 7).Value = 8.418600821238126
$ws.Cells.Item(36, 2).Value = 0.7287194209349384
$ws.Cells.Item(36, 3).Value = 0.3375848360084654
$ws.Cells.Item(36, 4).Value = 3.082599426703578
$ws.Cells.Item(36, 5).Value = 0.4998867070740569
$ws.Cells.Item(36, 7).Value = 4.64879039072104
$ws.Cells.Item(37, 2).Value = 1.505614041169197
$ws.Cells.Item(37, 3).Value = 1.65323645889881
$ws.Cells.Item(37, 4).Value = 3.082599426703578
$ws.Cells.Item(37, 5).Value = 6.48142807727062
$ws.Cells.Item(37, 7).Value = 12.7228780040422
$ws.Cells.Item(38, 2).Value = 3.182878228561681
$ws.Cells.Item(38, 3).Value = 1.65323645889881
$ws.Cells.Item(38, 4).Value = 0.1529057820181812
$ws.Cells.Item(38, 5).Value = 0.4998867070740569
$ws.Cells.Item(38, 7).Value = 5.488907176552729
$ws.Cells.Item(39, 2).Value = 3.182878228561681
$ws.Cells.Item(39, 3).Value = 1.65323645889881
$ws.Cells.Item(39, 4).Value = 16.98373111632243
$ws.Cells.Item(39, 5).Value = 0.4998867070740569
$ws.Cells.Item(39, 7).Value = 22.31973251085698
$ws.Cells.Item(40, 2).Value = 1.505614041169197
$ws.Cells.Item(40, 3).Value = 1.65323645889881
$ws.Cells.Item(40, 4).Value = 0.7127328510149897
$ws.Cells.Item(40, 5).Value = 0.4998867070740569
$ws.Cells.Item(40, 7).Value = 4.371470058157054
$ws.Cells.Item(41, 2).Value = 1.505614041169197
$ws.Cells.Item(41, 3).Value = 1.65323645889881
$ws.Cells.Item(41, 4).Value = 0.7127328510149897
$ws.Cells.Item(41, 5).Value = 0.4998867070740569
$ws.Cells.Item(41, 7).Value = 4.371470058157054
$ws.Cells.Item(42, 2).Value = 3.182878228561681
$ws.Cells.Item(42, 3).Value = 1.65323645889881
$ws.Cells.Item(42, 4).Value = 0.7127328510149897
$ws.Cells.Item(42, 5).Value = 0.4998867070740569
$ws.Cells.Item(42, 7).Value = 6.048734245549538
$ws.Cells.Item(43, 2).Value = 3.182878228561681
$ws.Cells.Item(43, 3).Value = 0.3375848360084654
$ws.Cells.Item(43, 4).Value = 0.1529057820181812
$ws.Cells.Item(43, 5).Value = 0.4998867070740569
$ws.Cells.Item(43, 7).Value = 4.173255553662385
$ws.Cells.Item(44, 2).Value = 3.182878228561681
$ws.Cells.Item(44, 3).Value = 1.65323645889881
$ws.Cells.Item(44, 4).Value = 16.98373111632243
$ws.Cells.Item(44, 5).Value = 0.4998867070740569
$ws.Cells.Item(44, 7).Value = 22.31973251085698
